# Actualización automática 2025-06-16 15:20:08
#
# Updates sales figures for "RIOS CARRION ANGEL BENIGNO" / group "PORCELANATO"
# for the month of "junio" across the three sheets of the workbook:
#   - VENTAS POR GRUPO     (per-client breakdown by product group)
#   - VENTA MENSUAL        (per-client breakdown by month)
#   - CUMPLIMIENTO MENSUAL (budget vs. sales roll-up)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: VENTAS POR GRUPO
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# CERAMICAS AL COSTO S.A.S. - PORCELANATO
$wsGrupo.Range("M6").Value = 15693.84

# ZAMBRANO ANGELA MARIA - PORCELANATO (was 0, now has sales)
$wsGrupo.Range("M18").Value = 4481.57

# Summary row: count of clients (of 17) with nonzero PORCELANATO sales
$wsGrupo.Range("M19").Value = "2 de 17"

# Column F width nudges from 13 -> 14 characters
$wsGrupo.Columns.Item(6).ColumnWidth = 13.1

# ---------------------------------------------------------------------
# Sheet: VENTA MENSUAL
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# CERAMICAS AL COSTO S.A.S. - junio
$wsMensual.Range("F6").Value = 15693.84

# ZAMBRANO ANGELA MARIA - junio (was 0, now has sales)
$wsMensual.Range("F18").Value = 4481.57

# TOTAL row - junio
$wsMensual.Range("F19").Value = 21419.92

# Column F width nudges from 13 -> 14 characters
$wsMensual.Columns.Item(6).ColumnWidth = 13.1

# ---------------------------------------------------------------------
# Sheet: CUMPLIMIENTO MENSUAL
# ---------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# PORCELANATO row (16): VENTA, POR CUMPLIR, CUMPLIMIENTO
$wsCumpl.Range("D16").Value = 20152.73
$wsCumpl.Range("E16").Value = 8057.110000000001
$wsCumpl.Range("F16").Value = 0.7143865402994133

# TOTAL row (19): VENTA, POR CUMPLIR, CUMPLIMIENTO
$wsCumpl.Range("D19").Value = 21419.92
$wsCumpl.Range("E19").Value = 25799.38386304603
$wsCumpl.Range("F19").Value = 0.4536263402384314

# Column width nudges: D 13 -> 14, F 26 -> 24
$wsCumpl.Columns.Item(4).ColumnWidth = 13.1
$wsCumpl.Columns.Item(6).ColumnWidth = 23.1
